$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: AngleSharpParser passing rate 0.8125 -> 1, remove "Incorrect encoding detection" problem
$ws.Range("C2").Value = 1
$ws.Range("D2").ClearContents()

# Row 3: CsQueryParserTest passing rate 0.75 -> 0.93, problem text changes
$ws.Range("C3").Value = 0.93
$ws.Range("D3").Value = "Slow"

# Row 4: HtmlAgilityPackParser passing rate 0.375 -> 1, remove "Not working" problem
$ws.Range("C4").Value = 1
$ws.Range("D4").ClearContents()

# Update sheet view: scroll so A10 is top-left, and selection is D18
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D18").Select()
